# "Always cite newest BT" -- update the `source` column so that every row
# cites the newest matching IQB-Bildungstrend report instead of the report
# that was contemporaneous with that row's year.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$koeller2010 = "Köller, O., Knigge, M., & Tesch, B. (2010). *Sprachliche Kompetenzen im Ländervergleich.* Waxmann. http://www.iqb.hu-berlin.de/bt/LV08_09/LV_ZF_0809c.pdf"
$bt2018math  = "Stanat, P., Schipolowski, S., Mahler, N., Weirich, S. & Henschel, S. (2019). *IQB-Bildungstrend 2018. Mathematische und naturwissenschaftliche Kompetenzen am Ende der Sekundarstufe I im zweiten Ländervergleich.* Waxmann Verlag. https://directory.doabooks.org/handle/20.500.12854/50672 "
$bt2021both  = "Stanat, P., Schipolowski, S., Schneider, R., Sachse, K. A., Weirich, S. & Henschel, S. (Hrsg.). (2022). *IQB-Bildungstrend 2021: Kompetenzen in den Fächern Deutsch und Mathematik am Ende der 4. Jahrgangsstufe im dritten Ländervergleich.* Waxmann Verlag. https://directory.doabooks.org/handle/20.500.12854/94704 "
$bt2022lang  = "Stanat, P., Schipolowski, S., Schneider, R., Weirich, S., Henschel, S. & Sachse, K. A. (Hrsg.). (2023). *IQB-Bildungstrend 2022: Sprachliche Kompetenzen am Ende der 9. Jahrgangsstufe im dritten Ländervergleich.* Waxmann. https://elibrary.utb.de/doi/book/10.31244/9783830997771 "

$ws.Range("B2").Value = $koeller2010
$ws.Range("B3").Value = $bt2021both
$ws.Range("B4").Value = $bt2018math
$ws.Range("B5").Value = $bt2022lang
$ws.Range("B6").Value = $bt2021both
$ws.Range("B7").Value = $bt2018math
$ws.Range("B8").Value = $bt2021both
$ws.Range("B9").Value = $bt2022lang

# Remove the now-unused/obsolete shared-string entries by dropping the
# references no longer used (Excel prunes unreferenced shared strings when
# it rewrites sharedStrings.xml on save).

$ws.Range("B10").Select
$ws.Application.ActiveWindow.Zoom = 87
